$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login Credentials")

function Set-TextValue($cell, $text) {
    # Force the literal string into the cell without Excel's "smart" numeric
    # coercion (comma-grouped digit strings like ",,,000052964," would
    # otherwise be parsed as a number), then restore the default cell style.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 33: update last name / full name / order-number tracker cells
$ws.Cells.Item(33, 4).Value = "Sethi D D D"
$ws.Cells.Item(33, 5).Value = "Nalini Sethi D D D"
Set-TextValue $ws.Cells.Item(33, 6) ",,,000052964,"

# Row 34: new credentials record (Vidhur Chopra)
$ws.Cells.Item(34, 1).Value = "unzoq@gmail.com"
$ws.Cells.Item(34, 2).Value = "3Eflpr385@"
$ws.Cells.Item(34, 3).Value = "Vidhur"
$ws.Cells.Item(34, 4).Value = "Chopra"
$ws.Cells.Item(34, 5).Value = "Vidhur Chopra"

# Row 35: new credentials record (Urmila Talwar)
$ws.Cells.Item(35, 1).Value = "fsjol@gmail.com"
$ws.Cells.Item(35, 2).Value = "PE9vzx758^"
$ws.Cells.Item(35, 3).Value = "Urmila"
$ws.Cells.Item(35, 4).Value = "Talwar D"
$ws.Cells.Item(35, 5).Value = "Urmila Talwar D"
$ws.Cells.Item(35, 6).Value = ","
